$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 539.65717
$ws.Range("J17").Value = 539.65717
$ws.Range("L17").Value = 1618.97151
$ws.Range("N17").Value = -1954.97151
$ws.Range("H29").Value = 1316.6
$ws.Range("I29").Value = 27.666666
$ws.Range("J29").Value = 3250
$ws.Range("K29").Value = 82.99999800000001
$ws.Range("L29").Value = 9750
$ws.Range("M29").Value = 198.000002
$ws.Range("N29").Value = -10312
$ws.Range("H38").Value = 2335.9412
$ws.Range("I38").Value = 193.90909
$ws.Range("J38").Value = 6263
$ws.Range("K38").Value = 581.72727
$ws.Range("L38").Value = 18789
$ws.Range("M38").Value = -209.72727
$ws.Range("N38").Value = -19533
$ws.Range("H70").Value = 1798
$ws.Range("J70").Value = 1755.1666
$ws.Range("L70").Value = 5265.4998
$ws.Range("N70").Value = -5805.4998
$ws.Range("H73").Value = 1798
$ws.Range("J73").Value = 1755.1666
$ws.Range("L73").Value = 5265.4998
$ws.Range("N73").Value = -7137.4998
$ws.Range("H106").Value = 27502912
$ws.Range("I106").Value = 33848636
$ws.Range("K106").Value = 33848636
$ws.Range("M106").Value = -33848005
$ws.Range("H137").Value = 3054.25
$ws.Range("I137").Value = 2727.5715
$ws.Range("J137").Value = 4034.2856
$ws.Range("K137").Value = 8182.7145
$ws.Range("L137").Value = 12102.8568
$ws.Range("M137").Value = -5632.7145
$ws.Range("N137").Value = -17202.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3292.0833
$ws.Range("I2").Value = 3050.6
$ws.Range("K2").Value = 3050.6
$ws.Range("M2").Value = -2937.6
$ws.Range("H32").Value = 7286.197
$ws.Range("I32").Value = 3508.3035
$ws.Range("K32").Value = 3508.3035
$ws.Range("M32").Value = -3221.3035
$ws.Range("H74").Value = 6349
$ws.Range("I74").Value = 2401.2354
$ws.Range("K74").Value = 2401.2354
$ws.Range("M74").Value = -1527.2354
$ws.Range("H77").Value = 6349
$ws.Range("I77").Value = 2401.2354
$ws.Range("K77").Value = 12006.177
$ws.Range("M77").Value = -7638.177
$ws.Range("H88").Value = 824.3684
$ws.Range("I88").Value = 612
$ws.Range("K88").Value = 612
$ws.Range("M88").Value = -206
$ws.Range("H91").Value = 824.3684
$ws.Range("I91").Value = 612
$ws.Range("K91").Value = 612
$ws.Range("M91").Value = 792
$ws.Range("H110").Value = 20834578
$ws.Range("I110").Value = 1260.9524
$ws.Range("J110").Value = 166667790
$ws.Range("K110").Value = 1260.9524
$ws.Range("L110").Value = 166667790
$ws.Range("M110").Value = 784.0476000000001
$ws.Range("N110").Value = -166671880
$ws.Range("H116").Value = 3292.0833
$ws.Range("I116").Value = 3050.6
$ws.Range("K116").Value = 3050.6
$ws.Range("M116").Value = -756.5999999999999
$ws.Range("H122").Value = 2661.1765
$ws.Range("I122").Value = 2856.5
$ws.Range("K122").Value = 8569.5
$ws.Range("M122").Value = -6119.5
$ws.Range("H132").Value = 2287.535
$ws.Range("I132").Value = 1268.9395
$ws.Range("K132").Value = 3806.8185
$ws.Range("M132").Value = -1276.8185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3292.0833
$ws.Range("I3").Value = 3050.6
$ws.Range("K3").Value = 3050.6
$ws.Range("M3").Value = -2936.6
$ws.Range("H107").Value = 31250948
$ws.Range("I107").Value = 752.46155
$ws.Range("K107").Value = 752.46155
$ws.Range("M107").Value = 1167.53845
$ws.Range("H134").Value = 2202.074
$ws.Range("I134").Value = 1321.909
$ws.Range("K134").Value = 3965.727
$ws.Range("M134").Value = -1430.727
$ws.Range("H140").Value = 131249.75
$ws.Range("J140").Value = 131249.75
$ws.Range("L140").Value = 131249.75
$ws.Range("N140").Value = -141609.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2074
$ws.Range("I16").Value = 1098.6666
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1098.6666
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -811.6666
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 10744.647
$ws.Range("I31").Value = 4310.6
$ws.Range("J31").Value = 59000
$ws.Range("K31").Value = 4310.6
$ws.Range("L31").Value = 59000
$ws.Range("M31").Value = -4015.6
$ws.Range("N31").Value = -59590
$ws.Range("H34").Value = 10744.647
$ws.Range("I34").Value = 4310.6
$ws.Range("J34").Value = 59000
$ws.Range("K34").Value = 4310.6
$ws.Range("L34").Value = 59000
$ws.Range("M34").Value = -4108.6
$ws.Range("N34").Value = -59404
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H113").Value = 2074
$ws.Range("I113").Value = 1098.6666
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1098.6666
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1071.3334
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 2661.0908
$ws.Range("I132").Value = 1996.1818
$ws.Range("J132").Value = 3990.9092
$ws.Range("K132").Value = 5988.5454
$ws.Range("L132").Value = 11972.7276
$ws.Range("M132").Value = -3458.5454
$ws.Range("N132").Value = -17032.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 225.4
$ws.Range("I7").Value = 32.333332
$ws.Range("J7").Value = 515
$ws.Range("K7").Value = 96.999996
$ws.Range("L7").Value = 1545
$ws.Range("M7").Value = 15.000004
$ws.Range("N7").Value = -1769
$ws.Range("H82").Value = 18010.5
$ws.Range("I82").Value = 18010.5
$ws.Range("K82").Value = 54031.5
$ws.Range("M82").Value = -53625.5
$ws.Range("H85").Value = 18010.5
$ws.Range("I85").Value = 18010.5
$ws.Range("K85").Value = 54031.5
$ws.Range("M85").Value = -52627.5
$ws.Range("H97").Value = 1401
$ws.Range("J97").Value = 568
$ws.Range("L97").Value = 1704
$ws.Range("N97").Value = -2696
$ws.Range("H122").Value = 970.8
$ws.Range("I122").Value = 966.3333
$ws.Range("J122").Value = 977.5
$ws.Range("K122").Value = 8696.9997
$ws.Range("L122").Value = 8797.5
$ws.Range("M122").Value = -6246.9997
$ws.Range("N122").Value = -13697.5
$ws.Range("H131").Value = 650642.7
$ws.Range("I131").Value = 811.82355
$ws.Range("J131").Value = 1654926.9
$ws.Range("K131").Value = 2435.47065
$ws.Range("L131").Value = 4964780.699999999
$ws.Range("M131").Value = 2604.52935
$ws.Range("N131").Value = -4974860.699999999
$ws.Range("H140").Value = 1573.4615
$ws.Range("J140").Value = 1583.3334
$ws.Range("L140").Value = 4750.0002
$ws.Range("N140").Value = -15110.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13119.5
$ws.Range("I80").Value = 12799.5
$ws.Range("K80").Value = 12799.5
$ws.Range("M80").Value = -11801.5
$ws.Range("H83").Value = 13119.5
$ws.Range("I83").Value = 12799.5
$ws.Range("K83").Value = 63997.5
$ws.Range("M83").Value = -59005.5
$ws.Range("H113").Value = 2567.6956
$ws.Range("I113").Value = 2563.4285
$ws.Range("K113").Value = 2563.4285
$ws.Range("M113").Value = -393.4285
$ws.Range("H131").Value = 79998
$ws.Range("J131").Value = 79998
$ws.Range("L131").Value = 79998
$ws.Range("N131").Value = -90078
$ws.Range("H138").Value = 99996.5
$ws.Range("J138").Value = 99996.5
$ws.Range("L138").Value = 99996.5
$ws.Range("N138").Value = -110276.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1939.8
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 2666.3333
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 2666.3333
$ws.Range("M22").Value = -555
$ws.Range("N22").Value = -3256.3333
$ws.Range("H27").Value = 1939.8
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 2666.3333
$ws.Range("K27").Value = 850
$ws.Range("L27").Value = 2666.3333
$ws.Range("M27").Value = -743
$ws.Range("N27").Value = -2880.3333
$ws.Range("H40").Value = 2889.5151
$ws.Range("I40").Value = 2178.5
$ws.Range("J40").Value = 9999.666999999999
$ws.Range("K40").Value = 2178.5
$ws.Range("L40").Value = 9999.666999999999
$ws.Range("M40").Value = -2042.5
$ws.Range("N40").Value = -10271.667
$ws.Range("H61").Value = 5034.933
$ws.Range("J61").Value = 10002.25
$ws.Range("L61").Value = 10002.25
$ws.Range("N61").Value = -10406.25
$ws.Range("H113").Value = 5034.933
$ws.Range("J113").Value = 10002.25
$ws.Range("L113").Value = 10002.25
$ws.Range("N113").Value = -14342.25
$ws.Range("H122").Value = 7186.1816
$ws.Range("I122").Value = 5506.125
$ws.Range("K122").Value = 16518.375
$ws.Range("M122").Value = -14068.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 199950
$ws.Range("J103").Value = 199950
$ws.Range("L103").Value = 199950
$ws.Range("N103").Value = -202294
$ws.Range("H105").Value = 119824.8
$ws.Range("J105").Value = 119824.8
$ws.Range("L105").Value = 119824.8
$ws.Range("N105").Value = -126812.8
$ws.Range("H113").Value = 1438.4286
$ws.Range("I113").Value = 325.64285
$ws.Range("K113").Value = 976.9285500000001
$ws.Range("M113").Value = 1193.07145
$ws.Range("H122").Value = 3939.25
$ws.Range("I122").Value = 3135.7856
$ws.Range("K122").Value = 9407.356800000001
$ws.Range("M122").Value = -6957.356800000001
